$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the English header values with their "verbeuse" French equivalents.
$ws.Range("G2").Value = "Restaurant avec cuisine sur place"
$ws.Range("H2").Value = "Concédée"
$ws.Range("I2").Value = "Public"

# type_production (G2) gets a Times New Roman font.
$ws.Range("G2").Font.Name = "Times New Roman"

# type_gestion (H2) and modele_economique (I2) switch to General format with wrapped text.
$ws.Range("H2").NumberFormat = "General"
$ws.Range("H2").WrapText = $true

$ws.Range("I2").NumberFormat = "General"
$ws.Range("I2").WrapText = $true

# Move the active selection to I2.
$null = $ws.Range("I2").Select()
